$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows 28-32: Key / Locator pairs appended to the Locators sheet
$ws.Range("A28").Value = "SuccessMessage_Xpath"
$ws.Range("B28").Value = "//ul[@class='messages']/li/ul/li"

$ws.Range("A29").Value = "HomePage_Account_Xpath"
$ws.Range("B29").Value = "//*[@id='lnkAccount']/a"

$ws.Range("A30").Value = "HomePage_Login_Xpath"
$ws.Range("B30").Value = "//*[@id='divAccount']/ul/li[1]/a"

$ws.Range("A31").Value = "No_Of_Customer_Reviews_Xpath"
$ws.Range("B31").Value = "(//div[@class='col-lg-9 customer_reviews pull-right']/div)"

$ws.Range("A32").Value = "Review_Date_Xpath"
$ws.Range("B32").Value = "(//li[@class='rvw_title block clear']/div)[2]"

# B29 and B30 pick up the highlighted (blue Courier New) locator-value
# formatting already used elsewhere in the sheet (e.g. B13/B14).
$ws.Range("B13").Copy()
$ws.Range("B29").PasteSpecial(-4122)
$ws.Range("B29").Value = "//*[@id='lnkAccount']/a"

$ws.Range("B13").Copy()
$ws.Range("B30").PasteSpecial(-4122)
$ws.Range("B30").Value = "//*[@id='divAccount']/ul/li[1]/a"

# Match the workbook's final selection / active cell
$ws.Range("A32").Select()
